# "Generate Report for Handback"
# The handback round-trip completed for both xliff files (66936f4b... and
# 9353dc2d...), in both target languages (zh-cn, de-de). Refresh the
# generated status report:
#   * flip the Status column from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is used (Overview's
#     per-language status columns, and each language tab's Status column),
#   * record the new "Latest Target File" / "Latest Handback File" names on
#     each language tab (with github links matching the existing
#     "Source File Name" links),
#   * stamp the de-de tab's "Latest Handback DateTime" with the handback
#     timestamp,
#   * widen the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$handbackTime = "2016-11-14 07:28:37"

$md1 = "66936f4b-add1-4a72-89ed-a1f62cca8d11.md"
$md2 = "9353dc2d-e6dd-41c1-9eb3-a31a8bb5e4ee.md"
$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8edf9892d6af96aa4f796d572767bed47a1ced36/e2e/66936f4b-add1-4a72-89ed-a1f62cca8d11.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8edf9892d6af96aa4f796d572767bed47a1ced36/e2e/9353dc2d-e6dd-41c1-9eb3-a31a8bb5e4ee.md"

# width (in "characters") that renders as the file's raw <col width> of 40 /
# 30 in this engine's column-width quantization (xml_width = (round(6*cw)+5)/6)
$wFileName = 38.333333333333336
$wStatus   = 29.166666666666668

# ---------------------------------------------------------------------
# Overview sheet: widen the per-language Status columns (E = zh-cn,
# F = de-de) and refresh their text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $wStatus
$wsOverview.Columns.Item(6).ColumnWidth = $wStatus

# ---------------------------------------------------------------------
# Per-language tabs (zh-cn, de-de): Status column C, Latest Target File
# column I, Latest Handback File column J, Latest Handback DateTime
# column K.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Range("C2:C3").Value = $newStatus
    $ws.Columns.Item(3).ColumnWidth = $wStatus
    $ws.Columns.Item(9).ColumnWidth = $wFileName
    $ws.Columns.Item(10).ColumnWidth = $wFileName
}

# zh-cn: row 2 -> 66936f4b..., row 3 -> 9353dc2d...
$wsZhCn.Range("I2").Value = $md1
$wsZhCn.Range("I2").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $url1, "", "", $md1)
$wsZhCn.Range("J2").Value = "66936f4b-add1-4a72-89ed-a1f62cca8d11.6177965fef03a409ee5e1abd4f43afcbbee97989.zh-cn.xlf"

$wsZhCn.Range("I3").Value = $md2
$wsZhCn.Range("I3").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $url2, "", "", $md2)
$wsZhCn.Range("J3").Value = "9353dc2d-e6dd-41c1-9eb3-a31a8bb5e4ee.88b463b0c21f7923e09fa69092b6a989b5f31d2c.zh-cn.xlf"

# de-de: row 2 -> 66936f4b..., row 3 -> 9353dc2d..., plus the handback
# timestamp now that the round trip is done.
$wsDeDe.Range("I2").Value = $md1
$wsDeDe.Range("I2").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $url1, "", "", $md1)
$wsDeDe.Range("J2").Value = "66936f4b-add1-4a72-89ed-a1f62cca8d11.6177965fef03a409ee5e1abd4f43afcbbee97989.de-de.xlf"
$wsDeDe.Range("K2").Value = $handbackTime

$wsDeDe.Range("I3").Value = $md2
$wsDeDe.Range("I3").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $url2, "", "", $md2)
$wsDeDe.Range("J3").Value = "9353dc2d-e6dd-41c1-9eb3-a31a8bb5e4ee.88b463b0c21f7923e09fa69092b6a989b5f31d2c.de-de.xlf"
$wsDeDe.Range("K3").Value = $handbackTime
